$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.776.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.390.98'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '504.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.96%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.552'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.398.37'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("E10").Value = '  +1.06%  '
$ws.Range("E11").Value = '  -0.54%  '
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.68'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.815.56'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.682.81'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.70'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.29%  '
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.409.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '309.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.72%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.13'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.17%  '
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.378'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.151'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '176.05'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0726'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.61%  '
$ws.Range("E32").Value = '  -0.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.13'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.58%  '
$ws.Range("E34").Value = '  -4.73%  '
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.20'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.82'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.826'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.85'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.82%  '
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '131.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.39'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.87'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.569'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '250.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0910'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.03%  '
$ws.Range("E49").Value = '  -0.93%  '
$ws.Range("E50").Value = '  +0.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.85%  '
